$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

$ws.Range("D2").Value = 0.0003334479406476021
$ws.Range("E2").Value = 0.1232835752889514
$ws.Range("G2").Value = 0.006893382407724857
$ws.Range("H2").Value = 0.01055803522467613
$ws.Range("I2").Value = 0.06034007668495178
$ws.Range("J2").Value = 0.03968131495639682
$ws.Range("K2").Value = 0.001510561909526587
$ws.Range("D3").Value = 0.01014292286708951
$ws.Range("E3").Value = 0.1390881040133536
$ws.Range("G3").Value = 0.00676274998113513
$ws.Range("H3").Value = 0.02089785737916827
$ws.Range("I3").Value = 0.05615905625745654
$ws.Range("J3").Value = 0.0493462230078876
$ws.Range("K3").Value = 0.001542714424431324
$ws.Range("D4").Value = 0.009386209305375814
$ws.Range("E4").Value = 0.1433697752654552
$ws.Range("G4").Value = 0.006645115558058023
$ws.Range("H4").Value = 0.0192387979477644
$ws.Range("I4").Value = 0.06180515419691801
$ws.Range("J4").Value = 0.05014666821807623
$ws.Range("K4").Value = 0.001486039720475674
$ws.Range("D5").Value = 0.0004838951863348484
$ws.Range("E5").Value = 0.1254716287367046
$ws.Range("G5").Value = 0.006760122254490852
$ws.Range("H5").Value = 0.01092662895098329
$ws.Range("I5").Value = 0.0610449630767107
$ws.Range("J5").Value = 0.04063326586037874
$ws.Range("K5").Value = 0.001608534716069698
$ws.Range("D6").Value = 0.01102607138454914
$ws.Range("E6").Value = 1.400171426124871
$ws.Range("G6").Value = 0.01685901638120413
$ws.Range("H6").Value = 0.04502810444682837
$ws.Range("I6").Value = 1.242255682125688
$ws.Range("J6").Value = 0.07556350808590651
$ws.Range("K6").Value = 0.006039711181074381
$ws.Range("D8").Value = 0.0003334479406476021
$ws.Range("E8").Value = 0.1232835752889514
$ws.Range("G8").Value = 0.006893382407724857
$ws.Range("H8").Value = 0.01055803522467613
$ws.Range("I8").Value = 0.06034007668495178
$ws.Range("J8").Value = 0.03968131495639682
$ws.Range("K8").Value = 0.001510561909526587
$ws.Range("D9").Value = 0.01014292286708951
$ws.Range("E9").Value = 0.1390881040133536
$ws.Range("G9").Value = 0.00676274998113513
$ws.Range("H9").Value = 0.02089785737916827
$ws.Range("I9").Value = 0.05615905625745654
$ws.Range("J9").Value = 0.0493462230078876
$ws.Range("K9").Value = 0.001542714424431324
$ws.Range("D10").Value = 0.009386209305375814
$ws.Range("E10").Value = 0.1433697752654552
$ws.Range("G10").Value = 0.006645115558058023
$ws.Range("H10").Value = 0.0192387979477644
$ws.Range("I10").Value = 0.06180515419691801
$ws.Range("J10").Value = 0.05014666821807623
$ws.Range("K10").Value = 0.001486039720475674
$ws.Range("D11").Value = 0.0004838951863348484
$ws.Range("E11").Value = 0.1254716287367046
$ws.Range("G11").Value = 0.006760122254490852
$ws.Range("H11").Value = 0.01092662895098329
$ws.Range("I11").Value = 0.0610449630767107
$ws.Range("J11").Value = 0.04063326586037874
$ws.Range("K11").Value = 0.001608534716069698
$ws.Range("D12").Value = 0.01102607138454914
$ws.Range("E12").Value = 1.400171426124871
$ws.Range("G12").Value = 0.01685901638120413
$ws.Range("H12").Value = 0.04502810444682837
$ws.Range("I12").Value = 1.242255682125688
$ws.Range("J12").Value = 0.07556350808590651
$ws.Range("K12").Value = 0.006039711181074381
